# Denmark Superligaen workbook update
# The source data had several match rows whose underlying "id" (column B)
# and all associated match data (columns B..AC, i.e. everything except the
# running index in column A) got re-synced to the correct record - in
# practice this manifests as neighbouring rows being swapped (two rows
# exchange all of their B:AC content) and, in one case, three rows being
# rotated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# Simple two-row swaps (column A, the running index, is left untouched)
Swap-Rows 27 28
Swap-Rows 32 33
Swap-Rows 80 81
Swap-Rows 110 112
Swap-Rows 207 208
Swap-Rows 261 262
Swap-Rows 352 353
Swap-Rows 448 449
Swap-Rows 454 455
Swap-Rows 490 491
Swap-Rows 499 500
Swap-Rows 539 540

# Three-row rotation: row305 <- old row306, row306 <- old row307, row307 <- old row305
$range305 = $ws.Range("B305:AC305")
$range306 = $ws.Range("B306:AC306")
$range307 = $ws.Range("B307:AC307")

$vals305 = $range305.Value2
$vals306 = $range306.Value2
$vals307 = $range307.Value2

$range305.Value2 = $vals306
$range306.Value2 = $vals307
$range307.Value2 = $vals305

Write-Host "Row swaps/rotation applied"
